$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("E2").Value = 22.86000000000013
$ws.Range("G2").Value = 0.0008672570368000176
$ws.Range("H2").Value = 0.003189769579054117

# I2 goes from a numeric value to an explicit empty text cell (matches J2's
# existing empty inline/shared string). A plain Value = "" clears the cell
# entirely in this engine, so force "text" typing via the classic leading
# apostrophe trick, then restore the default style so no stray formatting
# (quote-prefix) sticks around.
$ws.Range("I2").Value = "'"
$ws.Range("I2").Style = "Normal"

$ws.Range("K2").Value = 4.285738813768346
$ws.Range("L2").Value = "[1.8616150154943636, 6.709862612042329]"
$ws.Range("M2").Value = 0.0005630062990584772
$ws.Range("N2").Value = 0.0005630062990584772
$ws.Range("O2").Value = -1.622684493746079
$ws.Range("P2").Value = "[-2.402579366709311, -0.8427896207828471]"
$ws.Range("Q2").Value = 5.162665002478306 / 100000
$ws.Range("R2").Value = 0.0001032533000495661
$ws.Range("S2").Value = 13.51368051370517
$ws.Range("T2").Value = "[11.936923485523801, 15.090437541886542]"
$ws.Range("W2").Value = 5.903783783783815
$ws.Range("X2").Value = 3.066306306306321
$ws.Range("Y2").Value = 8.741261261261311

# --- Row 3 ---
$ws.Range("B3").Value = 0
$ws.Range("E3").Value = 23.23000000000019
$ws.Range("G3").Value = 2.135243587497726 / 10000000
$ws.Range("H3").Value = 65.14664309817849 / 10000000

$ws.Range("I3").Value = "'"
$ws.Range("I3").Style = "Normal"

$ws.Range("K3").Value = 5.743253057176473
$ws.Range("L3").Value = "[3.1642328040892718, 8.322273310263673]"
$ws.Range("M3").Value = 1.487533883381964 / 100000
$ws.Range("N3").Value = 2.975067766763928 / 100000
$ws.Range("O3").Value = -0.4402632347373086
$ws.Range("P3").Value = "[-0.9308422677303092, 0.0503157982556921]"
$ws.Range("Q3").Value = 0.0784704680022128
$ws.Range("R3").Value = 0.0784704680022128
$ws.Range("S3").Value = 13.12437842712644
$ws.Range("T3").Value = "[11.706533048286435, 14.542223805966447]"
$ws.Range("W3").Value = 1.627727727727745
$ws.Range("X3").Value = -0.1860260260260236
$ws.Range("Y3").Value = 3.441481481481513
